$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B-G (N column is G)
$data = @(
    @{ row = 2;  B = 0.08696016040226752; C = 0.2244626349384274;  D = 0.06217355075105529; E = 0.2493462467153963; F = 0.24189316101004;    G = 15 },
    @{ row = 3;  B = 0.3031818451638543;  C = 0.3412803276417196;  D = 0.1783459041844743;  E = 0.4223101990059846; F = 0.3050817571715199;  G = 14 },
    @{ row = 4;  B = 0.4726833070249565;  C = 0.4978299971988559;  D = 0.4530240298790372;  E = 0.6730705979903128; F = 0.4987257742493065;  G = 13 },
    @{ row = 5;  B = 0.6055399787264154;  C = 0.6086332724655883;  D = 0.5906225627394526;  E = 0.7685197217635034; F = 0.4942695954492832;  G = 12 },
    @{ row = 6;  B = 0.5882471822302764;  C = 0.5882471822302764;  D = 0.4329084706557275;  E = 0.6579578030966177; F = 0.3091295773284307;  G = 11 },
    @{ row = 7;  B = 0.4671934630864089;  C = 0.4696202125542849;  D = 0.2831907857527659;  E = 0.5321567304401644; F = 0.2685786741804888;  G = 10 },
    @{ row = 8;  B = 0.3964777679756588;  C = 0.4142371153650239;  D = 0.2111306574106052;  E = 0.4594895618081059; F = 0.2463291325150163;  G = 9 },
    @{ row = 9;  B = 0.430140519527954;   C = 0.430140519527954;   D = 0.2113118878782925;  E = 0.4596867279771002; F = 0.1776210167919809;  G = 6 },
    @{ row = 10; B = 0.3316630433627284;  C = 0.3316630433627284;  D = 0.1188298149235608;  E = 0.344717007012362;  F = 0.1150832780485532;  G = 3 },
    @{ row = 11; B = 0.6206497229122814;  C = 0.6206497229122814;  D = 0.3852060785510916;  E = 0.6206497229122814; F = $null;               G = 1 }
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    if ($item.F -eq $null) {
        $ws.Range("F$r").ClearContents()
    } else {
        $ws.Range("F$r").Value = $item.F
    }
    $ws.Range("G$r").Value = $item.G
}
